# Apply the cryptos-list refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.396.30"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "3.910.78"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.02"
$ws.Range("E5").Value = "  +9.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.77"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000334"
$ws.Range("E11").Value = "  -4.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.20"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "4.534.47"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.27"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("D15").Value = "3.921.55"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.06"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.22"
$ws.Range("E17").Value = "  +8.38%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.135"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.76"
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("D20").Value = "69.407.69"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.16"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("E22").Value = "  -4.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.22"
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.51"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.05"
$ws.Range("E25").Value = "  +10.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.54"
$ws.Range("E26").Value = "  -3.24%  "
$ws.Range("E27").Value = "  -3.29%  "
$ws.Range("E28").Value = "  -2.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "685.15"
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.17"
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "67.78"
$ws.Range("E33").Value = "  +11.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.444"
$ws.Range("E34").Value = "  +12.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.02"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.04"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0850"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0482"
$ws.Range("E41").Value = "  -3.14%  "
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.11"
$ws.Range("E43").Value = "  +3.99%  "
$ws.Range("E44").Value = "  -5.02%  "
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.141"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0357"
$ws.Range("E47").Value = "  +10.84%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.02"
$ws.Range("E48").Value = "  +7.96%  "
$ws.Range("D49").Value = "2.748.06"
$ws.Range("E49").Value = "  +13.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.10"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  -2.71%  "
